# Daniels assesments merged in to my structure
$wb = $excel.ActiveWorkbook

# Sheet named "Peer  and self assessment" holds the two criterion tables
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# Criterion 1 Online collaboration - Daniel's row (row 7)
$ws.Range("B7").Value = "Excellent"
$ws.Range("C7").Value = "Admin on discord and very active to support as well as very active in communications."

# Criterion 1 International Collaboration - Daniel's row (row 20)
$ws.Range("B20").Value = "Excellent"
$ws.Range("C20").Value = "Active collaborator, motivated"
